# Generate Report for Handback
# Applies the "Latest Target File" hyperlink + new handback datetime / error
# detail for the a8801518-... row (row 8) on both the zh-cn and de-de sheets,
# and widens the "Error Detail" column (P) to fit the new long message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ff204ea289dfe3dfb43b9657bc72cda2f53dd17/e2e/a8801518-fe82-443d-950e-d58457ad583b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d699a95aba506606051165c9842c891e46c0765b/e2e/a8801518-fe82-443d-950e-d58457ad583b.md."

$targetMdDisplay = "a8801518-fe82-443d-950e-d58457ad583b.md"
$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d699a95aba506606051165c9842c891e46c0765b/e2e/a8801518-fe82-443d-950e-d58457ad583b.md"

function Update-HandbackRow {
    param(
        $ws,
        [string]$handbackFile,
        [string]$handbackDateTime
    )

    # Widen the "Error Detail" column (P) so the long message is readable.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I8 = "Latest Target File" -> becomes a hyperlink to the handback .md,
    # matching the same pattern already used elsewhere in the sheet (e.g. I2).
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetMdUrl, "", "", $targetMdDisplay)

    # J8 = "Latest Handback File"
    $ws.Range("J8").Value = $handbackFile

    # K8 = "Latest Handback DateTime"
    $ws.Range("K8").Value = $handbackDateTime

    # P8 = "Error Detail"
    $ws.Range("P8").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow -ws $wsZhCn `
    -handbackFile "a8801518-fe82-443d-950e-d58457ad583b.cf3f21f53b6960ea7fc5a91d240fae5055d4be82.zh-cn.xlf" `
    -handbackDateTime "2016-08-19 20:48:08"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow -ws $wsDeDe `
    -handbackFile "a8801518-fe82-443d-950e-d58457ad583b.cf3f21f53b6960ea7fc5a91d240fae5055d4be82.de-de.xlf" `
    -handbackDateTime "2016-08-19 20:48:15"
